$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1/J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style (bold, centered, bordered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J67
$data = @{
    2 = @(8, 8)
    3 = @(5, 5)
    4 = @(7, 7)
    5 = @(8, 8)
    6 = @(9, 9)
    7 = @(7, 7)
    8 = @(8, 8)
    9 = @(7, 7)
    10 = @(6, 6)
    11 = @(6, 7)
    12 = @(1, 3)
    13 = @(7, 7)
    14 = @(7, 7)
    15 = @(7, 8)
    16 = @(7, 7)
    17 = @(8, 8)
    18 = @(8, 8)
    19 = @(6, 6)
    20 = @(7, 8)
    21 = @(7, 7)
    22 = @(8, 8)
    23 = @(9, 9)
    24 = @(6, 6)
    25 = @(8, 8)
    26 = @(8, 8)
    27 = @(8, 8)
    28 = @(8, 8)
    29 = @(8, 8)
    30 = @(8, 8)
    31 = @(6, 6)
    32 = @(7, 7)
    33 = @(5, 6)
    34 = @(9, 9)
    35 = @(7, 7)
    36 = @(8, 8)
    37 = @(8, 8)
    38 = @(7, 7)
    39 = @(9, 9)
    40 = @(8, 8)
    41 = @(8, 8)
    42 = @(6, 6)
    43 = @(8, 8)
    44 = @(8, 8)
    45 = @(10, 10)
    46 = @(8, 8)
    47 = @(8, 8)
    48 = @(8, 8)
    49 = @(6, 6)
    50 = @(7, 8)
    51 = @(6, 6)
    52 = @(8, 8)
    53 = @(8, 8)
    54 = @(7, 7)
    55 = @(6, 7)
    56 = @(7, 7)
    57 = @(7, 7)
    58 = @(9, 9)
    59 = @(8, 8)
    60 = @(8, 8)
    61 = @(8, 8)
    62 = @(6, 6)
    63 = @(5, 5)
    64 = @(8, 8)
    65 = @(7, 7)
    66 = @(4, 4)
    67 = @(5, 5)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item([int]$r, 9).Value = $vals[0]
    $ws.Cells.Item([int]$r, 10).Value = $vals[1]
}

Write-Host "done"